$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the date-block (current rows 305:306),
# pushing the existing rows 305:326 down to 307:328.
$ws.Rows("305:306").Insert()

# Copy the (now shifted) old row 305/306 data - which landed on 307/308 -
# into the newly-inserted blank rows 305/306, so all the static columns
# (A,B,C,E,F,G,H,I,N,O,Q,R) are populated consistently with the rest of
# the block. We'll overwrite D,J,K,L,M,P afterwards with the new values.
$ws.Range("A307:R308").Copy()
$ws.Range("A305").PasteSpecial()

# New "Primera" record (row 305)
$ws.Range("D305").Value = 44746
$ws.Range("J305").Value = 1200
$ws.Range("K305").Value = 450
$ws.Range("L305").Value = 500
$ws.Range("M305").Value = 475
$ws.Range("P305").Value = 119

# New "Segunda" record (row 306)
$ws.Range("D306").Value = 44746
$ws.Range("J306").Value = 1200
$ws.Range("K306").Value = 450
$ws.Range("L306").Value = 500
$ws.Range("M306").Value = 475
$ws.Range("P306").Value = 95
